$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 176.6
$ws.Cells.Item(12, 9).Value = 146
$ws.Cells.Item(12, 10).Value = 299
$ws.Cells.Item(12, 11).Value = 146
$ws.Cells.Item(12, 12).Value = 299
$ws.Cells.Item(12, 13).Value = 24
$ws.Cells.Item(12, 14).Value = -639
$ws.Cells.Item(17, 8).Value = 1828.8572
$ws.Cells.Item(17, 10).Value = 1828.8572
$ws.Cells.Item(17, 12).Value = 5486.571599999999
$ws.Cells.Item(17, 14).Value = -5822.571599999999
$ws.Cells.Item(28, 8).Value = 2854.625
$ws.Cells.Item(28, 9).Value = 2854.625
$ws.Cells.Item(28, 11).Value = 2854.625
$ws.Cells.Item(28, 13).Value = -2369.625
$ws.Cells.Item(40, 8).Value = 94501.03
$ws.Cells.Item(40, 9).Value = 377657.75
$ws.Cells.Item(40, 10).Value = 3890.88
$ws.Cells.Item(40, 11).Value = 377657.75
$ws.Cells.Item(40, 12).Value = 3890.88
$ws.Cells.Item(40, 13).Value = -377482.75
$ws.Cells.Item(40, 14).Value = -4240.88
$ws.Cells.Item(43, 8).Value = 1100
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).Value = $null
$ws.Cells.Item(74, 8).Value = 5548.9
$ws.Cells.Item(74, 9).Value = 5548.9
$ws.Cells.Item(74, 11).Value = 5548.9
$ws.Cells.Item(74, 13).Value = -4612.9
$ws.Cells.Item(77, 8).Value = 5548.9
$ws.Cells.Item(77, 9).Value = 5548.9
$ws.Cells.Item(77, 11).Value = 27744.5
$ws.Cells.Item(77, 13).Value = -23064.5
$ws.Cells.Item(111, 8).Value = 2262.5715
$ws.Cells.Item(111, 9).Value = 2915.2
$ws.Cells.Item(111, 10).Value = 631
$ws.Cells.Item(111, 11).Value = 8745.599999999999
$ws.Cells.Item(111, 12).Value = 1893
$ws.Cells.Item(111, 13).Value = -5678.599999999999
$ws.Cells.Item(111, 14).Value = -8027
$ws.Cells.Item(138, 8).Value = 2842.5967
$ws.Cells.Item(138, 9).Value = 1725.0769
$ws.Cells.Item(138, 10).Value = 3649.6943
$ws.Cells.Item(138, 11).Value = 5175.2307
$ws.Cells.Item(138, 12).Value = 10949.0829
$ws.Cells.Item(138, 13).Value = -35.23070000000007
$ws.Cells.Item(138, 14).Value = -21229.0829

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2259.4912
$ws.Cells.Item(32, 9).Value = 1515.3269
$ws.Cells.Item(32, 11).Value = 1515.3269
$ws.Cells.Item(32, 13).Value = -1228.3269
$ws.Cells.Item(61, 8).Value = 2911.0417
$ws.Cells.Item(61, 9).Value = 2588.8096
$ws.Cells.Item(61, 11).Value = 2588.8096
$ws.Cells.Item(61, 13).Value = -2376.8096
$ws.Cells.Item(63, 8).Value = 3977.2222
$ws.Cells.Item(63, 9).Value = 2500
$ws.Cells.Item(63, 11).Value = 2500
$ws.Cells.Item(63, 13).Value = -1814
$ws.Cells.Item(66, 8).Value = 3977.2222
$ws.Cells.Item(66, 9).Value = 2500
$ws.Cells.Item(66, 11).Value = 12500
$ws.Cells.Item(66, 13).Value = -9068
$ws.Cells.Item(74, 8).Value = 2005.7778
$ws.Cells.Item(74, 9).Value = 1970.8235
$ws.Cells.Item(74, 11).Value = 1970.8235
$ws.Cells.Item(74, 13).Value = -1096.8235
$ws.Cells.Item(77, 8).Value = 2005.7778
$ws.Cells.Item(77, 9).Value = 1970.8235
$ws.Cells.Item(77, 11).Value = 9854.1175
$ws.Cells.Item(77, 13).Value = -5486.1175
$ws.Cells.Item(110, 8).Value = 2412.4666
$ws.Cells.Item(110, 9).Value = 2168.2307
$ws.Cells.Item(110, 11).Value = 2168.2307
$ws.Cells.Item(110, 13).Value = -123.2307000000001
$ws.Cells.Item(136, 8).Value = 2911.0417
$ws.Cells.Item(136, 9).Value = 2588.8096
$ws.Cells.Item(136, 11).Value = 7766.4288
$ws.Cells.Item(136, 13).Value = -5216.4288

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 13336962
$ws.Cells.Item(134, 9).Value = 1354.25
$ws.Cells.Item(134, 11).Value = 4062.75
$ws.Cells.Item(134, 13).Value = -1527.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3795.0688
$ws.Cells.Item(31, 9).Value = 2154.7
$ws.Cells.Item(31, 10).Value = 7440.3335
$ws.Cells.Item(31, 11).Value = 2154.7
$ws.Cells.Item(31, 12).Value = 7440.3335
$ws.Cells.Item(31, 13).Value = -1859.7
$ws.Cells.Item(31, 14).Value = -8030.3335
$ws.Cells.Item(34, 8).Value = 3795.0688
$ws.Cells.Item(34, 9).Value = 2154.7
$ws.Cells.Item(34, 10).Value = 7440.3335
$ws.Cells.Item(34, 11).Value = 2154.7
$ws.Cells.Item(34, 12).Value = 7440.3335
$ws.Cells.Item(34, 13).Value = -1952.7
$ws.Cells.Item(34, 14).Value = -7844.3335
$ws.Cells.Item(58, 8).Value = 2765.7837
$ws.Cells.Item(58, 9).Value = 2399.7727
$ws.Cells.Item(58, 10).Value = 3302.6
$ws.Cells.Item(58, 11).Value = 2399.7727
$ws.Cells.Item(58, 12).Value = 3302.6
$ws.Cells.Item(58, 13).Value = -2196.7727
$ws.Cells.Item(58, 14).Value = -3708.6
$ws.Cells.Item(99, 8).Value = 1819.4
$ws.Cells.Item(99, 9).Value = 1819.4
$ws.Cells.Item(99, 11).Value = 1819.4
$ws.Cells.Item(99, 13).Value = -321.4000000000001
$ws.Cells.Item(126, 8).Value = 1819.4
$ws.Cells.Item(126, 9).Value = 1819.4
$ws.Cells.Item(126, 11).Value = 5458.200000000001
$ws.Cells.Item(126, 13).Value = -2988.200000000001
$ws.Cells.Item(132, 8).Value = 1428.0857
$ws.Cells.Item(132, 9).Value = 1126.2963
$ws.Cells.Item(132, 11).Value = 3378.8889
$ws.Cells.Item(132, 13).Value = -848.8888999999999
$ws.Cells.Item(136, 8).Value = 2765.7837
$ws.Cells.Item(136, 9).Value = 2399.7727
$ws.Cells.Item(136, 10).Value = 3302.6
$ws.Cells.Item(136, 11).Value = 7199.3181
$ws.Cells.Item(136, 12).Value = 9907.799999999999
$ws.Cells.Item(136, 13).Value = -4649.3181
$ws.Cells.Item(136, 14).Value = -15007.8
$ws.Cells.Item(140, 8).Value = 557499.5
$ws.Cells.Item(140, 10).Value = 557499.5
$ws.Cells.Item(140, 12).Value = 557499.5
$ws.Cells.Item(140, 14).Value = -567859.5
$ws.Cells.Item(141, 8).Value = 639409.6
$ws.Cells.Item(141, 10).Value = 673185.25
$ws.Cells.Item(141, 12).Value = 673185.25
$ws.Cells.Item(141, 14).Value = -683545.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 749
$ws.Cells.Item(68, 9).Value = 749
$ws.Cells.Item(68, 11).Value = 2247
$ws.Cells.Item(68, 13).Value = -1436
$ws.Cells.Item(71, 8).Value = 749
$ws.Cells.Item(71, 9).Value = 749
$ws.Cells.Item(71, 11).Value = 6741
$ws.Cells.Item(71, 13).Value = -2685
$ws.Cells.Item(80, 8).Value = 2490
$ws.Cells.Item(80, 9).Value = 2490
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 7470
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = -6534
$ws.Cells.Item(80, 14).Value = $null
$ws.Cells.Item(83, 8).Value = 2490
$ws.Cells.Item(83, 9).Value = 2490
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 22410
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = -17730
$ws.Cells.Item(83, 14).Value = $null
$ws.Cells.Item(121, 8).Value = 951.8
$ws.Cells.Item(121, 9).Value = 253.83333
$ws.Cells.Item(121, 11).Value = 761.49999
$ws.Cells.Item(121, 13).Value = 548.50001
$ws.Cells.Item(132, 8).Value = 1668
$ws.Cells.Item(132, 9).Value = 933.4
$ws.Cells.Item(132, 10).Value = 2586.25
$ws.Cells.Item(132, 11).Value = 8400.6
$ws.Cells.Item(132, 12).Value = 23276.25
$ws.Cells.Item(132, 13).Value = -5870.6
$ws.Cells.Item(132, 14).Value = -28336.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 15963.4
$ws.Cells.Item(70, 9).Value = 36840.332
$ws.Cells.Item(70, 10).Value = 4220.125
$ws.Cells.Item(70, 11).Value = 36840.332
$ws.Cells.Item(70, 12).Value = 4220.125
$ws.Cells.Item(70, 13).Value = -36570.332
$ws.Cells.Item(70, 14).Value = -4760.125
$ws.Cells.Item(73, 8).Value = 15963.4
$ws.Cells.Item(73, 9).Value = 36840.332
$ws.Cells.Item(73, 10).Value = 4220.125
$ws.Cells.Item(73, 11).Value = 36840.332
$ws.Cells.Item(73, 12).Value = 4220.125
$ws.Cells.Item(73, 13).Value = -35904.332
$ws.Cells.Item(73, 14).Value = -6092.125
$ws.Cells.Item(80, 8).Value = 2340.4167
$ws.Cells.Item(80, 10).Value = 2958
$ws.Cells.Item(80, 12).Value = 2958
$ws.Cells.Item(80, 14).Value = -4954
$ws.Cells.Item(83, 8).Value = 2340.4167
$ws.Cells.Item(83, 10).Value = 2958
$ws.Cells.Item(83, 12).Value = 14790
$ws.Cells.Item(83, 14).Value = -24774
$ws.Cells.Item(107, 8).Value = 1763.7727
$ws.Cells.Item(107, 9).Value = 1819.5714
$ws.Cells.Item(107, 10).Value = 1666.125
$ws.Cells.Item(107, 11).Value = 1819.5714
$ws.Cells.Item(107, 12).Value = 1666.125
$ws.Cells.Item(107, 13).Value = 100.4286
$ws.Cells.Item(107, 14).Value = -5506.125
$ws.Cells.Item(126, 8).Value = 2363.6667
$ws.Cells.Item(126, 9).Value = 2214.9092
$ws.Cells.Item(126, 11).Value = 6644.7276
$ws.Cells.Item(126, 13).Value = -4174.7276
$ws.Cells.Item(132, 8).Value = 1726.2222
$ws.Cells.Item(132, 9).Value = 1733.6471
$ws.Cells.Item(132, 11).Value = 5200.9413
$ws.Cells.Item(132, 13).Value = -2670.9413

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1467
$ws.Cells.Item(22, 9).Value = 934.5
$ws.Cells.Item(22, 11).Value = 934.5
$ws.Cells.Item(22, 13).Value = -639.5
$ws.Cells.Item(27, 8).Value = 1467
$ws.Cells.Item(27, 9).Value = 934.5
$ws.Cells.Item(27, 11).Value = 934.5
$ws.Cells.Item(27, 13).Value = -827.5
$ws.Cells.Item(61, 8).Value = 1738.6364
$ws.Cells.Item(61, 9).Value = 1738.6364
$ws.Cells.Item(61, 11).Value = 1738.6364
$ws.Cells.Item(61, 13).Value = -1536.6364
$ws.Cells.Item(113, 8).Value = 1738.6364
$ws.Cells.Item(113, 9).Value = 1738.6364
$ws.Cells.Item(113, 11).Value = 1738.6364
$ws.Cells.Item(113, 13).Value = 431.3635999999999
$ws.Cells.Item(128, 8).Value = 76698.42999999999
$ws.Cells.Item(128, 10).Value = 76698.42999999999
$ws.Cells.Item(128, 12).Value = 76698.42999999999
$ws.Cells.Item(128, 14).Value = -86658.42999999999
$ws.Cells.Item(136, 8).Value = 3879.6667
$ws.Cells.Item(136, 9).Value = 3782.75
$ws.Cells.Item(136, 10).Value = 4189.8
$ws.Cells.Item(136, 11).Value = 11348.25
$ws.Cells.Item(136, 12).Value = 12569.4
$ws.Cells.Item(136, 13).Value = -8798.25
$ws.Cells.Item(136, 14).Value = -17669.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 11168.154
$ws.Cells.Item(126, 9).Value = 12671
$ws.Cells.Item(126, 10).Value = 9414.833000000001
$ws.Cells.Item(126, 11).Value = 38013
$ws.Cells.Item(126, 12).Value = 28244.499
$ws.Cells.Item(126, 13).Value = -35543
$ws.Cells.Item(126, 14).Value = -33184.499
$ws.Cells.Item(132, 8).Value = 2332.2666
$ws.Cells.Item(132, 9).Value = 1966.3889
$ws.Cells.Item(132, 11).Value = 5899.1667
$ws.Cells.Item(132, 13).Value = -3369.1667
